$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "imageLink" column (O): header label plus the same image
# URL for both data rows.
$ws.Range("O1").Value = "imageLink"
$ws.Range("O2").Value = "http://dev.dailytexanonline.com/sites/default/files/images/2015/01/a%20different%20name.jpg"
$ws.Range("O3").Value = "http://dev.dailytexanonline.com/sites/default/files/images/2015/01/a%20different%20name.jpg"

# Scroll the window so column C is the left-most visible column, then
# leave the new O3 cell as the active selection.
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 1
$ws.Range("O3").Select()
